$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.658.23"
$ws.Range("E2").Value = "  +0.37%  "
$ws.Range("D3").Value = "'1.879.82"
$ws.Range("E3").Value = "  -0.47%  "
$ws.Range("D4").Value = "'1.013"
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").Value = "'316.21"
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").Value = "'1.010"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("D7").Value = "'0.5106"
$ws.Range("E7").Value = "  -1.19%  "
$ws.Range("D8").Value = "'0.3931"
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.08406"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("E10").Value = "  -1.39%  "
$ws.Range("D11").Value = "'41.88"
$ws.Range("E11").Value = "  -0.63%  "
$ws.Range("D12").Value = "'6.265"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "'1.882.91"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "'20.46"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").Value = "'7.268"
$ws.Range("E15").Value = "  -0.35%  "
$ws.Range("D16").Value = "'1.012"
$ws.Range("E16").Value = "  -0.54%  "
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("D18").Value = "'91.56"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'0.06721"
$ws.Range("E19").Value = "  -0.49%  "
$ws.Range("D20").Value = "'17.80"
$ws.Range("E20").Value = "  -0.38%  "
$ws.Range("D21").Value = "'1.009"
$ws.Range("E21").Value = "  -0.80%  "
$ws.Range("E22").Value = "  -1.23%  "
$ws.Range("D23").Value = "'28.678.64"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").Value = "'11.15"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("D25").Value = "'2.249"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "'2.095.04"
$ws.Range("E26").Value = "  -0.32%  "
$ws.Range("D27").Value = "'162.01"
$ws.Range("E27").Value = "  +0.42%  "
$ws.Range("D28").Value = "'20.75"
$ws.Range("E28").Value = "  -0.67%  "
$ws.Range("D29").Value = "'2.364"
$ws.Range("E29").Value = "  -3.36%  "
$ws.Range("E30").Value = "  -0.47%  "
$ws.Range("D31").Value = "'0.1054"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("D32").Value = "'1.052"
$ws.Range("E32").Value = "  -0.23%  "
$ws.Range("D33").Value = "'5.812"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("D34").Value = "'3.617"
$ws.Range("E34").Value = "  -0.82%  "
$ws.Range("D35").Value = "'0.02471"
$ws.Range("E35").Value = "  +0.39%  "
$ws.Range("D36").Value = "'0.06553"
$ws.Range("E36").Value = "  -0.36%  "
$ws.Range("D37").Value = "'0.2180"
$ws.Range("E37").Value = "  -1.17%  "
$ws.Range("D38").Value = "'8.923"
$ws.Range("E38").Value = "  -5.08%  "
$ws.Range("D39").Value = "'1.270"
$ws.Range("E39").Value = "  +1.01%  "
$ws.Range("D40").Value = "'1.204"
$ws.Range("E40").Value = "  +0.52%  "
$ws.Range("E41").Value = "  -0.82%  "
$ws.Range("D42").Value = "'5.058"
$ws.Range("E42").Value = "  +1.03%  "
$ws.Range("D43").Value = "'11.17"
$ws.Range("E43").Value = "  -0.58%  "
$ws.Range("D44").Value = "'1.009"
$ws.Range("E44").Value = "  -0.85%  "
$ws.Range("D45").Value = "'0.6053"
$ws.Range("E45").Value = "  -1.13%  "
$ws.Range("E46").Value = "  -1.51%  "
$ws.Range("D47").Value = "'3.701"
$ws.Range("E47").Value = "  -0.52%  "
$ws.Range("D48").Value = "'2.030"
$ws.Range("E48").Value = "  +0.65%  "
$ws.Range("D49").Value = "'1.219"
$ws.Range("E49").Value = "  -0.22%  "
$ws.Range("D50").Value = "'122.48"
$ws.Range("E50").Value = "  +0.15%  "
$ws.Range("D51").Value = "'1.187"
$ws.Range("E51").Value = "  -7.98%  "
